# Apply the refreshed cryptos list (GitHub Actions data-pull) to the
# "Price" (D) and "Volume(1h)" (E) columns. Both columns are stored as
# plain text in the workbook (t="inlineStr"), so every Price write is
# forced to text via NumberFormat "@" (then the style is reset back to
# "Normal" so we do not leave a stray cell style behind) to stop Excel
# from reinterpreting values such as "215.09" as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.946.55'
$ws.Range("E2").Value = '  -0.37%  '

# Row 3
$ws.Range("D3").Value = '1.673.80'
$ws.Range("E3").Value = '  +0.92%  '

# Row 4
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '

# Row 6
$ws.Range("E6").Value = '  +1.75%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0619'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.20%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.27'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.66%  '

# Row 11
$ws.Range("E11").Value = '  +0.72%  '

# Row 12
$ws.Range("D12").Value = '1.909.64'
$ws.Range("E12").Value = '  +0.94%  '

# Row 13
$ws.Range("D13").Value = '1.691.23'
$ws.Range("E13").Value = '  +2.02%  '

# Row 14
$ws.Range("E14").Value = '  -0.09%  '

# Row 15
$ws.Range("E15").Value = '  +0.64%  '

# Row 16
$ws.Range("E16").Value = '  -0.09%  '

# Row 17
$ws.Range("D17").Value = '26.948.19'
$ws.Range("E17").Value = '  -0.40%  '

# Row 18
$ws.Range("E18").Value = '  +4.87%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '235.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.80%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0735'
$ws.Range("E20").Value = '  -0.48%  '

# Row 21
$ws.Range("E21").Value = '  +0.14%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.23%  '

# Row 23
$ws.Range("E23").Value = '  -1.15%  '

# Row 24
$ws.Range("E24").Value = '  -2.04%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.37%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.02%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.112'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.29%  '

# Row 29
$ws.Range("E29").Value = '  +0.24%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0498'
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = '  +0.05%  '

# Row 32
$ws.Range("E32").Value = '  +0.50%  '

# Row 33
$ws.Range("D33").Value = '1.488.96'
$ws.Range("E33").Value = '  -4.14%  '

# Row 34
$ws.Range("E34").Value = '  +1.99%  '

# Row 35
$ws.Range("E35").Value = '  +3.20%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.36%  '

# Row 37
$ws.Range("E37").Value = '  +0.69%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.896'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.58%  '

# Row 39
$ws.Range("E39").Value = '  +0.56%  '

# Row 40
$ws.Range("E40").Value = '  +8.17%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.84'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.38%  '

# Row 42
$ws.Range("E42").Value = '  +0.16%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.31'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.00%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '67.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.49%  '

# Row 45
$ws.Range("D45").Value = '1.816.03'
$ws.Range("E45").Value = '  +0.66%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.777'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.15%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.62'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.48%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.53'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.07%  '

# Row 49
$ws.Range("E49").Value = '  +1.81%  '

# Row 50
$ws.Range("E50").Value = '  +0.25%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.64%  '
